$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows 12-41 of the "Espinaca" price table.
# Columns: RowNumber, Year, Month, Day (Fecha), Volumen(J), PrecioMinimo(K),
#          PrecioMaximo(L), PrecioPromedioPonderado(M), PrecioPorKg(P)
# Every other column (A,B,C,E,F,G,H,I,N,O,Q,R) is constant for every record
# in this table, so it only needs to be (re)written for the brand-new rows
# 39-41 that did not previously exist.
$rows = @(
    @(12, 2021, 8, 20, 50, 7500, 8000, 7800, 780),
    @(13, 2021, 5, 26, 50, 6000, 6500, 6300, 630),
    @(14, 2021, 6, 15, 50, 6000, 6500, 6300, 630),
    @(15, 2021, 1, 8, 80, 7000, 7500, 7188, 719),
    @(16, 2021, 6, 2, 60, 6000, 6500, 6250, 625),
    @(17, 2021, 8, 13, 100, 7000, 7500, 7250, 725),
    @(18, 2021, 3, 18, 80, 7000, 8000, 7500, 750),
    @(19, 2021, 1, 29, 60, 7500, 8000, 7750, 775),
    @(20, 2020, 12, 29, 100, 8000, 9000, 8500, 850),
    @(21, 2021, 7, 6, 60, 7500, 8000, 7750, 775),
    @(22, 2021, 5, 20, 60, 6000, 6500, 6250, 625),
    @(23, 2021, 2, 11, 100, 8000, 8500, 8250, 825),
    @(24, 2021, 5, 4, 50, 7000, 7500, 7200, 720),
    @(25, 2021, 8, 27, 100, 7000, 7500, 7250, 725),
    @(26, 2021, 7, 28, 80, 7500, 8000, 7688, 769),
    @(27, 2021, 5, 11, 60, 7000, 7500, 7250, 725),
    @(28, 2021, 1, 13, 80, 7500, 8000, 7688, 769),
    @(29, 2021, 2, 4, 70, 7500, 8000, 7714, 771),
    @(30, 2021, 4, 27, 60, 6000, 6500, 6250, 625),
    @(31, 2021, 6, 4, 50, 6000, 6500, 6300, 630),
    @(32, 2021, 7, 22, 60, 9000, 10000, 9500, 950),
    @(33, 2021, 4, 6, 50, 10000, 11000, 10600, 1060),
    @(34, 2021, 8, 25, 100, 7000, 7500, 7250, 725),
    @(35, 2020, 12, 23, 80, 8000, 8500, 8250, 825),
    @(36, 2021, 4, 30, 100, 6000, 6500, 6250, 625),
    @(37, 2020, 11, 25, 100, 9000, 9500, 9250, 925),
    @(38, 2021, 2, 3, 60, 9000, 10000, 9500, 950),
    @(39, 2021, 6, 8, 50, 6000, 6500, 6300, 630),
    @(40, 2021, 6, 29, 100, 6000, 6500, 6250, 625),
    @(41, 2021, 5, 18, 60, 6500, 7000, 6750, 675)
)

foreach ($row in $rows) {
    $r = $row[0]
    $year = $row[1]
    $month = $row[2]
    $day = $row[3]
    $volumen = $row[4]
    $precioMin = $row[5]
    $precioMax = $row[6]
    $precioProm = $row[7]
    $precioKg = $row[8]

    if ($r -ge 39) {
        # Brand new rows: fill in the constant columns too.
        $ws.Cells.Item($r, 1).Value = 11
        $ws.Cells.Item($r, 2).Value = "Vega Monumental Concepción"
        $ws.Cells.Item($r, 3).Value = "Bíobío"
        # Match the existing date-formatted cells' number format before
        # writing the value so no extra/duplicate style gets created.
        $ws.Cells.Item($r, 4).NumberFormat = $ws.Cells.Item(12, 4).NumberFormat
        $ws.Cells.Item($r, 5).Value = 8
        $ws.Cells.Item($r, 6).Value = 100112012
        $ws.Cells.Item($r, 7).Value = "Espinaca"
        $ws.Cells.Item($r, 8).Value = "Sin especificar"
        $ws.Cells.Item($r, 9).Value = "Primera"
        $ws.Cells.Item($r, 14).Value = "`$/cuna 10 kilos"
        $ws.Cells.Item($r, 15).Value = "Región Metropolitana"
        $ws.Cells.Item($r, 17).Value = 10
        $ws.Cells.Item($r, 18).Value = "Hortaliza"
    }

    $fecha = Get-Date -Year $year -Month $month -Day $day -Hour 0 -Minute 0 -Second 0
    $ws.Cells.Item($r, 4).Value = $fecha
    $ws.Cells.Item($r, 10).Value = $volumen
    $ws.Cells.Item($r, 11).Value = $precioMin
    $ws.Cells.Item($r, 12).Value = $precioMax
    $ws.Cells.Item($r, 13).Value = $precioProm
    $ws.Cells.Item($r, 16).Value = $precioKg
}
